$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Temps passé" (time spent) is no longer 0 for every task - fill in actuals
$ws.Range("E6").Value = 30
$ws.Range("E7").Value = 10
$ws.Range("E8").Value = 60
$ws.Range("E9").Value = 120

# "Avancement" (progress) should be expressed as a 0-100 percentage instead
# of a 0-1 fraction, so multiply the existing ratio formula by 100
$ws.Range("G6").Formula = "=((E6*100%)/D6)*100"
$ws.Range("G7").Formula = "=((E7*100%)/D7)*100"
$ws.Range("G8").Formula = "=((E8*100%)/D8)*100"
$ws.Range("G9").Formula = "=((E9*100%)/D9)*100"

# Fill in the little summary block in column J (Estimation / Temps passé /
# Reste à faire totals)
$ws.Range("J6").Formula = "=SUM(D6:D9)"
$ws.Range("J7").Formula = "=SUM(E6:E9)"
$ws.Range("J8").Formula = "=J6-J7"

# Reposition the view: scroll so column B is the leftmost visible column,
# then leave the selection on J9 (last cell of the table)
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("J9").Select()

# Nudge the saved window position down a bit, matching the author's session
$excel.ActiveWindow.Top = 1200
